# Append new daily COVID overview rows (2021-11-27 .. 2021-12-05) to the
# bottom of the "covid_totals" sheet, continuing the existing table layout:
#   A: date (text)          B: areaType        C: areaCode
#   D: areaName              E: cumCasesByPublishDate
#   F: newCasesByPublishDate G: newDeaths28DaysByPublishDate
#   H: cumDeaths28DaysByPublishDate

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$areaType = "overview"
$areaCode = "K02000001"
$areaName = "United Kingdom"

$newRows = @(
    @{ Row = 473; Date = "2021-11-27"; E = 10110408; F = 39567; G = 131; H = 144724 },
    @{ Row = 474; Date = "2021-11-28"; E = 10146915; F = 37681; G = 51;  H = 144775 },
    @{ Row = 475; Date = "2021-11-29"; E = 10189059; F = 42583; G = 35;  H = 144810 },
    @{ Row = 476; Date = "2021-11-30"; E = 10228772; F = 39716; G = 159; H = 144969 },
    @{ Row = 477; Date = "2021-12-01"; E = 10276007; F = 48374; G = 171; H = 145140 },
    @{ Row = 478; Date = "2021-12-02"; E = 10329074; F = 53945; G = 141; H = 145281 },
    @{ Row = 479; Date = "2021-12-03"; E = 10379647; F = 50584; G = 143; H = 145424 },
    @{ Row = 480; Date = "2021-12-04"; E = 10421104; F = 42848; G = 127; H = 145551 },
    @{ Row = 481; Date = "2021-12-05"; E = 10464389; F = 43992; G = 54;  H = 145605 }
)

# Format column A for the new rows as Text first so the date-looking
# strings ("2021-11-27", ...) are kept verbatim instead of being
# auto-converted into date serial numbers, matching every other row
# already in the sheet.
$ws.Range("A473:A481").NumberFormat = "@"

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Date
    $ws.Cells.Item($row, 2).Value = $areaType
    $ws.Cells.Item($row, 3).Value = $areaCode
    $ws.Cells.Item($row, 4).Value = $areaName
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
}
